$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.926.61"
$ws.Range("E2").Value = "  +0.39%  "
$ws.Range("D3").Value = "1.812.39"
$ws.Range("E3").Value = "  +1.69%  "
$ws.Range("E4").Value = "  -0.22%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "310.45"
$ws.Range("E5").Value = "  -0.07%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.002"
$ws.Range("E6").Value = "  -0.21%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4958"
$ws.Range("E7").Value = "  -2.97%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3872"
$ws.Range("E8").Value = "  +2.85%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.09713"
$ws.Range("E9").Value = "  +24.86%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.102"
$ws.Range("E10").Value = "  +1.54%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "40.99"
$ws.Range("E11").Value = "  -0.51%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.473"
$ws.Range("E12").Value = "  +4.43%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.52"
$ws.Range("E13").Value = "  +1.70%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.002"
$ws.Range("E14").Value = "  -0.20%  "
$ws.Range("D15").Value = "1.810.12"
$ws.Range("E15").Value = "  +1.49%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.298"
$ws.Range("E16").Value = "  +1.88%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001133"
$ws.Range("E17").Value = "  +5.87%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "92.61"
$ws.Range("E18").Value = "  +0.75%  "
$ws.Range("E19").Value = "  +1.09%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.001"
$ws.Range("E20").Value = "  -0.25%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.11"
$ws.Range("E21").Value = "  +0.70%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.926"
$ws.Range("E22").Value = "  +0.31%  "
$ws.Range("D23").Value = "27.975.19"
$ws.Range("E23").Value = "  +0.38%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.12"
$ws.Range("E24").Value = "  +1.70%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.246"
$ws.Range("E25").Value = "  +0.00%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "158.95"
$ws.Range("E26").Value = "  +0.28%  "
$ws.Range("D27").Value = "2.022.56"
$ws.Range("E27").Value = "  +1.79%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.58"
$ws.Range("E28").Value = "  +1.95%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.399"
$ws.Range("E29").Value = "  +1.97%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "127.33"
$ws.Range("E30").Value = "  +3.85%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.1059"
$ws.Range("E31").Value = "  -1.99%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.040"
$ws.Range("E32").Value = "  +0.87%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.577"
$ws.Range("E33").Value = "  +1.79%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.618"
$ws.Range("E34").Value = "  -0.25%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.06736"
$ws.Range("E35").Value = "  -4.71%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "8.972"
$ws.Range("E36").Value = "  +4.79%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02329"
$ws.Range("E37").Value = "  +1.11%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2139"
$ws.Range("E38").Value = "  +0.95%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.944"
$ws.Range("E39").Value = "  -1.21%  "
$ws.Range("E40").Value = "  -2.32%  "
$ws.Range("E41").Value = "  +1.52%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.001"
$ws.Range("E42").Value = "  -0.15%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.147"
$ws.Range("E43").Value = "  -0.08%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.13"
$ws.Range("E44").Value = "  +0.14%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5876"
$ws.Range("E45").Value = "  -1.62%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.695"
$ws.Range("E46").Value = "  -0.68%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.276"
$ws.Range("E47").Value = "  -3.68%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "122.67"
$ws.Range("E48").Value = "  -3.38%  "
$ws.Range("E49").Value = "  +2.00%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.177"
$ws.Range("E50").Value = "  -2.69%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06787"
$ws.Range("E51").Value = "  +1.08%  "
